$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# D-column values that look like plain decimal numbers must be forced to
# text (NumberFormat "@") before assignment, otherwise Excel auto-converts
# them to numeric values and silently drops significant trailing zeros
# (e.g. "593.00" -> 593, "1.00" -> 1), which would not match the source data.
$ws.Range("D2").Value = "67.891.84"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "2.531.56"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.00"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.57"
$ws.Range("E6").Value = "  +4.55%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").Value = "2.530.52"
$ws.Range("E9").Value = "  -1.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.141"
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.16"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.344"
$ws.Range("E13").Value = "  -3.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.85"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("D15").Value = "2.982.80"
$ws.Range("E15").Value = "  -2.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000178"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "67.562.33"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D18").Value = "2.527.39"
$ws.Range("E18").Value = "  -2.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.11"
$ws.Range("E19").Value = "  +4.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.48"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "359.81"
$ws.Range("E21").Value = "  +2.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.20"
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.67"
$ws.Range("E23").Value = "  +1.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.99"
$ws.Range("E24").Value = "  +3.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.29"
$ws.Range("E26").Value = "  +3.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "70.75"
$ws.Range("E27").Value = "  +1.93%  "
$ws.Range("D28").Value = "2.657.74"
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").Value = "0.0₃0988"
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "552.93"
$ws.Range("E31").Value = "  +4.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.31"
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("E33").Value = "  +1.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.86"
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.131"
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.48"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.15"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.78"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.59"
$ws.Range("E40").Value = "  +1.43%  "
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.355"
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.18"
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.55"
$ws.Range("E44").Value = "  +4.34%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "148.23"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.560"
$ws.Range("E47").Value = "  -1.37%  "
$ws.Range("D48").Value = "0.0₆0280"
$ws.Range("E48").Value = "  -2.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.71"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0759"
$ws.Range("E51").Value = "  -0.55%  "
